# Automated "Disponibilidad" refresh: each run re-stamps the "Ultimo" (Fecha)
# column for the most-recent availability-check block, and the timestamps
# that used to belong to the newer blocks cascade down to the next-older
# block (rows 2-15 -> 16-29 -> 30-43), matching the commit
# "Update automatico via Actualizar 02-22-2021 12-30-25".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newest block (rows 2-15): stamped with the new "now" timestamp.
$ws.Range("D2:D15").Value = 44249.52093606021

# Middle block (rows 16-29): takes on what used to be the newest block's
# timestamp.
$ws.Range("D16:D29").Value = 44249.49960440972

# Oldest block (rows 30-43): takes on what used to be the middle block's
# timestamp.
$ws.Range("D30:D43").Value = 44249.47827144676
